$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate the newly-added "Stadium Population" values in column D (rows 2-33).
# A handful of rows (matching the number format already used in column F for the
# same row) carry the thousands-separator number format; the rest stay General.
$ws.Range("D2").Value = 246709
$ws.Range("D2").NumberFormat = "#,##0"

$ws.Range("D3").Value = 486290
$ws.Range("D4").Value = 611648

$ws.Range("D5").Value = 29672
$ws.Range("D5").NumberFormat = "#,##0"

$ws.Range("D6").Value = 859035
$ws.Range("D7").Value = 29901
$ws.Range("D8").Value = 301301
$ws.Range("D9").Value = 385525

$ws.Range("D10").Value = 396394
$ws.Range("D10").NumberFormat = "#,##0"

$ws.Range("D11").Value = 704621
$ws.Range("D12").Value = 673104
$ws.Range("D13").Value = 105116
$ws.Range("D14").Value = 2312717
$ws.Range("D15").Value = 863002
$ws.Range("D16").Value = 892062
$ws.Range("D17").Value = 488943

$ws.Range("D18").Value = 113750
$ws.Range("D18").NumberFormat = "#,##0"

$ws.Range("D19").Value = 422331

$ws.Range("D20").Value = 17574
$ws.Range("D20").NumberFormat = "#,##0"

$ws.Range("D21").Value = 8622698

$ws.Range("D22").Value = 9928
$ws.Range("D22").NumberFormat = "#,##0"

$ws.Range("D23").Value = 9928
$ws.Range("D24").Value = 8622698
$ws.Range("D25").Value = 1580863
$ws.Range("D26").Value = 302407

$ws.Range("D27").Value = 3999759
$ws.Range("D27").NumberFormat = "#,##0"

$ws.Range("D28").Value = 92735
$ws.Range("D29").Value = 127134
$ws.Range("D30").Value = 724745
$ws.Range("D31").Value = 385430
$ws.Range("D32").Value = 667560
$ws.Range("D33").Value = 23078

# Row 24's "City Population" (column F) was also corrected to match column D.
$ws.Range("F24").Value = 8622698

# Update the active selection to reflect where the editor ended up after the edit.
$ws.Range("E34").Select()
